# Generate Report for Handoff
# Replace the old source-file UUID/name with the new one, update the
# generated xliff file names, and bump the handoff/handback timestamps.

$wb = $excel.ActiveWorkbook

$oldId = "b3b6aa9a-6094-4364-993f-2428770d7540"
$newId = "a923ed86-eea9-4a41-8447-fbda1f50c501"

$oldZhXlf = "$oldId.e41ba0e7646832f2efd8516a96e7d7ab891ec01a.zh-cn.xlf"
$newZhXlf = "$newId.e5862f00626cb73e5e628373a2c44d53a29d366b.zh-cn.xlf"

$oldDeXlf = "$oldId.e41ba0e7646832f2efd8516a96e7d7ab891ec01a.de-de.xlf"
$newDeXlf = "$newId.e5862f00626cb73e5e628373a2c44d53a29d366b.de-de.xlf"

$hyperlinkBaseUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/de896321e9c7639d46b5c34d0e32edcd39d3de5a/e2e/$oldId.md"

function Set-HyperlinkDisplay($ws, $cellAddr, $displayText) {
    $rng = $ws.Range($cellAddr)
    $rng.Hyperlinks.Delete()
    $hl = $rng.Hyperlinks.Item(1)
    $hl.Address = $hyperlinkBaseUrl
    $hl.TextToDisplay = $displayText
}

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# A2 - plain file name, no hyperlink on this sheet
$wsOverview.Range("A2").Value = "$newId.md"

# B2 - hyperlinked relative path
Set-HyperlinkDisplay $wsOverview "B2" "e2e\$newId.md"

# G2 - latest HO Xliff generate date
$wsOverview.Range("G2").Value = "2016-08-13 13:12:43"

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# A2 - hyperlinked file name
Set-HyperlinkDisplay $wsZhCn "A2" "$newId.md"

# G2 - latest handoff xliff file name
$wsZhCn.Range("G2").Value = $newZhXlf

# H2 - latest handoff datetime
$wsZhCn.Range("H2").Value = "2016-08-13 13:12:36"

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# A2 - hyperlinked file name
Set-HyperlinkDisplay $wsDeDe "A2" "$newId.md"

# G2 - latest handoff xliff file name
$wsDeDe.Range("G2").Value = $newDeXlf

# H2 - latest handback datetime
$wsDeDe.Range("H2").Value = "2016-08-13 13:12:43"
